# Generate Report for Handoff
# Adds a new localization-status row (f5d6eaa3-64cd-4225-a22c-3bb6154260ca.md)
# to all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileGuid   = "f5d6eaa3-64cd-4225-a22c-3bb6154260ca"
$fileName   = "$fileGuid.md"
$pathName   = "e2e\$fileGuid.md"
$zhXlf      = "$fileGuid.bdfbadb32e3c62dd61ea986c3e3f89136ab552db.zh-cn.xlf"
$deXlf      = "$fileGuid.bdfbadb32e3c62dd61ea986c3e3f89136ab552db.de-de.xlf"
$hoDateZh   = "2016-08-13 02:53:24"
$hoDateDe   = "2016-08-13 02:53:31"
$status     = "Ready for handoff"
$blobShaZh  = "97e12c54f4ae951ac545d3416db7b2bccefb5f60"
$blobShaDe  = "046510dc9e49bfcfbb22bbe63be40833b9b5cfe1"
$blobShaOv  = "bdfbadb32e3c62dd61ea986c3e3f89136ab552db"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ovWs = $wb.Worksheets.Item("Overview")
$ovTable = $ovWs.ListObjects.Item(1)
$ovRow = $ovTable.ListRows.Add()
$ovRange = $ovRow.Range

$ovRange.Cells.Item(1,1).Value = $fileName
$ovRange.Cells.Item(1,2).Value = $pathName
$ovRange.Cells.Item(1,3).Value = ".md"
$ovRange.Cells.Item(1,4).Value = ""
$ovRange.Cells.Item(1,5).Value = $status
$ovRange.Cells.Item(1,6).Value = $status
$ovRange.Cells.Item(1,7).Value = $hoDateDe
$ovRange.Cells.Item(1,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ovWs.Hyperlinks.Add($ovRange.Cells.Item(1,2), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$blobShaOv/e2e/$fileName", "", "", $pathName) | Out-Null

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zhWs = $wb.Worksheets.Item("zh-cn")
$zhTable = $zhWs.ListObjects.Item(1)
$zhRow = $zhTable.ListRows.Add()
$zhRange = $zhRow.Range

$zhRange.Cells.Item(1,1).Value  = $fileName
$zhRange.Cells.Item(1,2).Value  = ".md"
$zhRange.Cells.Item(1,3).Value  = $status
$zhRange.Cells.Item(1,4).Value  = "e2e"
$zhRange.Cells.Item(1,5).Value  = "ht"
$zhRange.Cells.Item(1,6).Value  = "False"
$zhRange.Cells.Item(1,7).Value  = $zhXlf
$zhRange.Cells.Item(1,8).Value  = $hoDateZh
$zhRange.Cells.Item(1,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhRange.Cells.Item(1,9).Value  = ""
$zhRange.Cells.Item(1,10).Value = ""
$zhRange.Cells.Item(1,11).Value = "0001-01-01 00:00:00"
$zhRange.Cells.Item(1,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhRange.Cells.Item(1,12).Value = ""
$zhRange.Cells.Item(1,13).Value = "True"
$zhRange.Cells.Item(1,14).Value = ""
$zhRange.Cells.Item(1,15).Value = "False"
$zhRange.Cells.Item(1,16).Value = ""

$zhWs.Hyperlinks.Add($zhRange.Cells.Item(1,1), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$blobShaZh/e2e/$fileName", "", "", $fileName) | Out-Null

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$deWs = $wb.Worksheets.Item("de-de")
$deTable = $deWs.ListObjects.Item(1)
$deRow = $deTable.ListRows.Add()
$deRange = $deRow.Range

$deRange.Cells.Item(1,1).Value  = $fileName
$deRange.Cells.Item(1,2).Value  = ".md"
$deRange.Cells.Item(1,3).Value  = $status
$deRange.Cells.Item(1,4).Value  = "e2e"
$deRange.Cells.Item(1,5).Value  = "ht"
$deRange.Cells.Item(1,6).Value  = "False"
$deRange.Cells.Item(1,7).Value  = $deXlf
$deRange.Cells.Item(1,8).Value  = $hoDateDe
$deRange.Cells.Item(1,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$deRange.Cells.Item(1,9).Value  = ""
$deRange.Cells.Item(1,10).Value = ""
$deRange.Cells.Item(1,11).Value = "0001-01-01 00:00:00"
$deRange.Cells.Item(1,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$deRange.Cells.Item(1,12).Value = ""
$deRange.Cells.Item(1,13).Value = "True"
$deRange.Cells.Item(1,14).Value = ""
$deRange.Cells.Item(1,15).Value = "False"
$deRange.Cells.Item(1,16).Value = ""

$deWs.Hyperlinks.Add($deRange.Cells.Item(1,1), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$blobShaDe/e2e/$fileName", "", "", $fileName) | Out-Null
